# Add a new row (row 3) to the active sheet containing the same shipping
# details as row 2, generated by Katalon AI.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new cells to be stored as text (matching the existing rows,
# which store every value - including numeric-looking ones - as text)
# before assigning values, then restore the default "Normal" style so we
# don't leave a custom number format applied to the row.
$ws.Range("A3:G3").NumberFormat = "@"

$ws.Range("A3").Value = "gf"
$ws.Range("B3").Value = "44"
$ws.Range("C3").Value = "fdg"
$ws.Range("D3").Value = "df"
$ws.Range("E3").Value = "5634653546546"
$ws.Range("F3").Value = "555"
$ws.Range("G3").Value = "fgg"

$ws.Range("A3:G3").Style = "Normal"
